# Rapise 6.6 note update
# - RVL sheet: collapse the "Map Range Data" param block down to just the
#   sheetName param (drop fromRow/fromCol/toRow/toCol), and rename the
#   Nav-related "Functions" actions to the new "Nav" object with shortened
#   action names (NavLaunch -> Launch, NavChangeCompany -> ChangeCompany,
#   NavNavigate -> Navigate, NavClose -> Close).
# - Cleanup sheet: drop two blank spacer rows near the top.

$wb = $excel.ActiveWorkbook

$wsRvl = $wb.Worksheets.Item("RVL")

# Remove the four now-unused Map/Range parameter rows (fromRow, fromCol,
# toRow, toCol). This shifts everything below up by 4 rows.
$wsRvl.Rows("9:12").Delete()

# Rename the "Functions" Nav actions to the "Nav" object with shorter
# action names (rows shifted up by 4 from their original 15-17, 22).
$wsRvl.Cells.Item(11, 3).Value = "Nav"
$wsRvl.Cells.Item(11, 4).Value = "Launch"

$wsRvl.Cells.Item(12, 3).Value = "Nav"
$wsRvl.Cells.Item(12, 4).Value = "ChangeCompany"

$wsRvl.Cells.Item(13, 3).Value = "Nav"
$wsRvl.Cells.Item(13, 4).Value = "Navigate"

$wsRvl.Cells.Item(18, 3).Value = "Nav"
$wsRvl.Cells.Item(18, 4).Value = "Close"

$wsCleanup = $wb.Worksheets.Item("Cleanup")

# Remove two blank spacer rows.
$wsCleanup.Rows("3:4").Delete()
